$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.185.39'
$ws.Range('D3').Value = '2.426.20'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.04'
$ws.Range('E5').Value = '  +1.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.01'
$ws.Range('E6').Value = '  +3.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.513'
$ws.Range('E7').Value = '  +1.23%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.501'
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.30'
$ws.Range('E10').Value = '  +2.86%  '
$ws.Range('E11').Value = '  +1.60%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.95'
$ws.Range('E12').Value = '  +3.54%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.124'
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.94'
$ws.Range('E14').Value = '  +1.84%  '
$ws.Range('D15').Value = '2.803.75'
$ws.Range('E15').Value = '  +1.87%  '
$ws.Range('D16').Value = '2.463.81'
$ws.Range('E16').Value = '  +4.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.835'
$ws.Range('E17').Value = '  +3.37%  '
$ws.Range('D18').Value = '44.130.82'
$ws.Range('E18').Value = '  +2.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.29'
$ws.Range('E19').Value = '  +1.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.42'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').Value = '0.0₃0905'
$ws.Range('E21').Value = '  +1.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.48'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.39'
$ws.Range('E23').Value = '  +2.08%  '
$ws.Range('E24').Value = '  +4.01%  '
$ws.Range('E25').Value = '  +1.26%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.18'
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.33'
$ws.Range('E28').Value = '  -1.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.61'
$ws.Range('E29').Value = '  +5.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.85'
$ws.Range('E30').Value = '  +4.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '18.67'
$ws.Range('E31').Value = '  +7.56%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.116'
$ws.Range('E32').Value = '  +10.34%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.18'
$ws.Range('E33').Value = '  +1.96%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  +2.35%  '
$ws.Range('E36').Value = '  +2.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.49'
$ws.Range('E37').Value = '  +4.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '129.08'
$ws.Range('E38').Value = '  +23.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.90'
$ws.Range('E39').Value = '  +3.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.30'
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('E41').Value = '  +0.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.19'
$ws.Range('E42').Value = '  -5.21%  '
$ws.Range('E43').Value = '  +2.27%  '
$ws.Range('D44').Value = '1.954.06'
$ws.Range('E44').Value = '  -0.38%  '
$ws.Range('E45').Value = '  +1.83%  '
$ws.Range('E46').Value = '  +4.34%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.39'
$ws.Range('E47').Value = '  +2.66%  '
$ws.Range('E48').Value = '  +9.40%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.41'
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.73'
$ws.Range('E50').Value = '  +2.43%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.16'
$ws.Range('E51').Value = '  +1.30%  '
